$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I values (LinearSVC benchmark) for rows 2-21
$values = @{
    2  = 0.53
    3  = 0.77
    4  = 0.62
    5  = 0.89
    6  = 0.98
    7  = 0.98
    8  = 0.65
    9  = 0.79
    10 = 0.83
    11 = 0.84
    12 = 1
    13 = 0.85
    14 = 0.96
    15 = 0.89
    16 = 0.71
    17 = 0.66
    18 = 0.74
    19 = 0.62
    20 = 0.86
    21 = 0.74
}

foreach ($row in ($values.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 9).Value = $values[$row]
}

# Update the selection to match the committed workbook state
$ws.Range("G25").Select()
